# Finish factorising the "one file one tab" series: the workbook already
# ends with tabs multiple_answers..multiple_answers14 (sheetId 7..21), each
# a copy of the same small "Identifiers" lookup table. Extend that series
# with 17 more tabs, multiple_answers15..multiple_answers31, by copying the
# current last sheet to the end of the workbook and renaming each copy.

$wb = $excel.ActiveWorkbook

for ($i = 15; $i -le 31; $i++) {
    $lastIndex = $wb.Worksheets.Count
    $lastSheet = $wb.Worksheets.Item($lastIndex)

    # Copy the last sheet, placing the new copy right after it (i.e. at the
    # new end of the workbook), matching Excel's "Move or Copy... > Create a
    # copy" behaviour used to extend this tab series.
    $lastSheet.Copy($null, $lastSheet)

    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = "multiple_answers$i"
}

# Restore the original active tab (copying shifts selection to the newly
# created sheet as a side effect).
$wb.Worksheets.Item(1).Activate()
